$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''86.283.11'
$ws.Range('E2').Value = '  +5.05%  '
$ws.Range('D3').Value = '''3.274.75'
$ws.Range('E3').Value = '  +3.37%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').Value = '''210.67'
$ws.Range('E5').Value = '  -3.10%  '
$ws.Range('D6').Value = '''625.52'
$ws.Range('E6').Value = '  +1.31%  '
$ws.Range('D7').Value = '''0.369'
$ws.Range('E7').Value = '  +28.23%  '
$ws.Range('D8').Value = '''0.653'
$ws.Range('E8').Value = '  +12.54%  '
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('D10').Value = '''3.270.71'
$ws.Range('E10').Value = '  +2.98%  '
$ws.Range('D11').Value = '''0.577'
$ws.Range('E11').Value = '  -2.22%  '
$ws.Range('D12').Value = '''0.177'
$ws.Range('E12').Value = '  +7.15%  '
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '''3.873.77'
$ws.Range('E14').Value = '  +2.89%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '''33.98'
$ws.Range('E15').Value = '  +6.52%  '
$ws.Range('D16').Value = '''5.28'
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').Value = '''86.055.21'
$ws.Range('E17').Value = '  +4.94%  '
$ws.Range('D18').Value = '''3.262.33'
$ws.Range('E18').Value = '  +2.55%  '
$ws.Range('D19').Value = '''14.10'
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('D20').Value = '''3.00'
$ws.Range('E20').Value = '  -7.32%  '
$ws.Range('D21').Value = '''429.29'
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').Value = '''8.93'
$ws.Range('E22').Value = '  +0.64%  '
$ws.Range('D23').Value = '''5.34'
$ws.Range('E23').Value = '  +4.55%  '
$ws.Range('D24').Value = '''7.17'
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('D25').Value = '''12.43'
$ws.Range('E25').Value = '  +4.54%  '
$ws.Range('D26').Value = '''5.10'
$ws.Range('E26').Value = '  -2.54%  '
$ws.Range('D27').Value = '''3.432.08'
$ws.Range('E27').Value = '  +2.95%  '
$ws.Range('D28').Value = '''75.84'
$ws.Range('E28').Value = '  -1.15%  '
$ws.Range('E29').Value = '  +6.85%  '
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').Value = '''0.175'
$ws.Range('E31').Value = '  +18.77%  '
$ws.Range('D32').Value = '''0.998'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('E33').Value = '  -2.09%  '
$ws.Range('D34').Value = '''545.38'
$ws.Range('E34').Value = '  -4.29%  '
$ws.Range('E35').Value = '  -3.86%  '
$ws.Range('E36').Value = '  -1.88%  '
$ws.Range('D37').Value = '''6.91'
$ws.Range('E37').Value = '  +11.53%  '
$ws.Range('D38').Value = '''0.137'
$ws.Range('E38').Value = '  -9.68%  '
$ws.Range('D39').Value = '''22.43'
$ws.Range('E39').Value = '  -0.58%  '
$ws.Range('D40').Value = '''0.998'
$ws.Range('E40').Value = '  -0.39%  '
$ws.Range('D41').Value = '''21.58'
$ws.Range('E41').Value = '  +3.56%  '
$ws.Range('D42').Value = '''0.393'
$ws.Range('E42').Value = '  -2.60%  '
$ws.Range('D43').Value = '''1.98'
$ws.Range('E43').Value = '  -1.01%  '
$ws.Range('E44').Value = '  -1.83%  '
$ws.Range('D45').Value = '''158.02'
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('D47').Value = '''178.80'
$ws.Range('E47').Value = '  -3.61%  '
$ws.Range('D48').Value = '''44.13'
$ws.Range('E48').Value = '  -1.16%  '
$ws.Range('D49').Value = '''1.30'
$ws.Range('E49').Value = '  -1.00%  '
$ws.Range('D50').Value = '''4.24'
$ws.Range('E50').Value = '  +1.76%  '
$ws.Range('D51').Value = '''0.621'
$ws.Range('E51').Value = '  -0.58%  '
